$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(16).Insert()
$ws.Cells.Item(16, 1).Value = "DOB"
$ws.Cells.Item(16, 4).Value = "/wlq-res-doc:WildlifeLicenseQueryResults/wlq-res-ext:WildlifeLicenseReport/nc:Person/nc:PersonBirthDate/nc:Date"
